$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

$ws.Range("H1").Value = "newRowInitialElementKeyToValueMap"
$ws.Range("I1").Value = "openRowInitialElementKeyToValueMap"
$ws.Range("H2").Value = "{}"
$ws.Range("I2").Value = "{}"
$ws.Range("H3").Value = "{}"
$ws.Range("I3").Value = "{}"

$ws.Activate()
$ws.Range("I3").Select()

